# "Documentation finis et ajustement appellation document"
#
# Fill in Responsable + Debut/Fin for every task row, add the remaining
# project tasks (countries + the voiture/css/html pages), and leave the
# selection where the author last left it.
#
# NOTE on ordering: the shared-string table is written/compacted in the
# order each distinct string is first assigned to a cell, so the statements
# below are deliberately ordered to reproduce that same first-use sequence
# (Page voiture, Hugo, Italie, Allemagne, Japon, Css Histoire,
# Css caracteristiques, Html Voitures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 : "Page voiture" (new task)
$ws.Range("A9").Value = "Page voiture"

# Rows 2-4 : existing tasks get a Responsable + dates
$ws.Range("B2").Value = "Hugo"
$ws.Range("C2").Value = 44935
$ws.Range("D2").Value = 44939

$ws.Range("B3").Value = "Hugo"
$ws.Range("C3").Value = 44935
$ws.Range("D3").Value = 44937

$ws.Range("B4").Value = "Hugo"
$ws.Range("C4").Value = 44935
$ws.Range("D4").Value = 44936

# Row 9 : dates + responsable for "Page voiture"
$ws.Range("B9").Value = "Hugo"
$ws.Range("C9").Value = 44937
$ws.Range("D9").Value = 44940

# Row 7 : Italie
$ws.Range("A7").Value = "Italie"
$ws.Range("B7").Value = "Hugo"
$ws.Range("C7").Value = 44936
$ws.Range("D7").Value = 44939

# Row 6 : Allemagne
$ws.Range("A6").Value = "Allemagne"
$ws.Range("B6").Value = "Hugo"
$ws.Range("C6").Value = 44936
$ws.Range("D6").Value = 44939

# Row 5 : Japon
$ws.Range("A5").Value = "Japon"
$ws.Range("B5").Value = "Hugo"
$ws.Range("C5").Value = 44936
$ws.Range("D5").Value = 44939

# Row 8 : "Page d'accueil" task keeps its label, gains responsable + dates
$ws.Range("B8").Value = "Hugo"
$ws.Range("C8").Value = 44935
$ws.Range("D8").Value = 44940

# Rows 10-12 : new tasks (Css Histoire / Css caracteristiques / Html Voitures)
$ws.Range("A10").Value = "Css Histoire"
$ws.Range("B10").Value = "Hugo"
$ws.Range("C10").Value = 44937
$ws.Range("D10").Value = 44940

$ws.Range("A11").Value = "Css caracteristiques"
$ws.Range("B11").Value = "Hugo"
$ws.Range("C11").Value = 44937
$ws.Range("D11").Value = 44940

$ws.Range("A12").Value = "Html Voitures"
$ws.Range("B12").Value = "Hugo"
$ws.Range("C12").Value = 44937
$ws.Range("D12").Value = 44940

# Rows 10-12 and the new trailing row 13 carry a (no-op) alignment tweak,
# which is what stamps the new 4th cell style used in the source workbook
$ws.Range("A10").WrapText = $false
$ws.Range("A11").WrapText = $false
$ws.Range("A12").WrapText = $false
$ws.Range("A13").WrapText = $false

# Leave the selection where the author left it when they saved
$ws.Range("E15").Select()
